# Auto-generated edit script applying the diff's cell-level value changes
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 844.25
$ws.Range("J17").Value = 844.25
$ws.Range("L17").Value = 2532.75
$ws.Range("N17").Value = -2868.75
$ws.Range("H19").Value = 1049.3846
$ws.Range("I19").Value = 844
$ws.Range("K19").Value = 844
$ws.Range("M19").Value = -669
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").Value = $null
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").Value = $null
$ws.Range("H132").Value = 1817.8334
$ws.Range("I132").Value = 1727.7059
$ws.Range("K132").Value = 5183.1177
$ws.Range("M132").Value = -2653.1177
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").Value = $null
$ws.Range("H141").Value = 5774.647
$ws.Range("I141").Value = 5611.3335
$ws.Range("J141").Value = 6999.5
$ws.Range("K141").Value = 16834.0005
$ws.Range("L141").Value = 20998.5
$ws.Range("M141").Value = -11654.0005
$ws.Range("N141").Value = -31358.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1403.6666
$ws.Range("I2").Value = 1355.5
$ws.Range("K2").Value = 1355.5
$ws.Range("M2").Value = -1242.5
$ws.Range("H45").Value = 2077
$ws.Range("I45").Value = 2077
$ws.Range("K45").Value = 2077
$ws.Range("M45").Value = -1700
$ws.Range("H61").Value = 7096.8
$ws.Range("I61").Value = 7096.8
$ws.Range("K61").Value = 7096.8
$ws.Range("M61").Value = -6884.8
$ws.Range("H74").Value = 1613
$ws.Range("I74").Value = 1517.8334
$ws.Range("K74").Value = 1517.8334
$ws.Range("M74").Value = -643.8334
$ws.Range("H77").Value = 1613
$ws.Range("I77").Value = 1517.8334
$ws.Range("K77").Value = 7589.166999999999
$ws.Range("M77").Value = -3221.166999999999
$ws.Range("H97").Value = 939.25
$ws.Range("I97").Value = 352.69232
$ws.Range("J97").Value = 3481
$ws.Range("K97").Value = 352.69232
$ws.Range("L97").Value = 3481
$ws.Range("M97").Value = 143.30768
$ws.Range("N97").Value = -4473
$ws.Range("H116").Value = 1403.6666
$ws.Range("I116").Value = 1355.5
$ws.Range("K116").Value = 1355.5
$ws.Range("M116").Value = 938.5
$ws.Range("H135").Value = 150000
$ws.Range("J135").Value = 150000
$ws.Range("L135").Value = 150000
$ws.Range("N135").Value = -160140
$ws.Range("H136").Value = 7096.8
$ws.Range("I136").Value = 7096.8
$ws.Range("K136").Value = 21290.4
$ws.Range("M136").Value = -18740.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1403.6666
$ws.Range("I3").Value = 1355.5
$ws.Range("K3").Value = 1355.5
$ws.Range("M3").Value = -1241.5
$ws.Range("H20").Value = 4319.6
$ws.Range("I20").Value = 4599.3335
$ws.Range("J20").Value = 3900
$ws.Range("K20").Value = 4599.3335
$ws.Range("L20").Value = 3900
$ws.Range("M20").Value = -4352.3335
$ws.Range("N20").Value = -4394
$ws.Range("H64").Value = 1299.4
$ws.Range("J64").Value = 1500
$ws.Range("L64").Value = 1500
$ws.Range("N64").Value = -1950
$ws.Range("H67").Value = 1299.4
$ws.Range("J67").Value = 1500
$ws.Range("L67").Value = 1500
$ws.Range("N67").Value = -3060
$ws.Range("H86").Value = 3143.7
$ws.Range("I86").Value = 3247.75
$ws.Range("J86").Value = 2727.5
$ws.Range("K86").Value = 3247.75
$ws.Range("L86").Value = 2727.5
$ws.Range("M86").Value = -2124.75
$ws.Range("N86").Value = -4973.5
$ws.Range("H89").Value = 3143.7
$ws.Range("I89").Value = 3247.75
$ws.Range("J89").Value = 2727.5
$ws.Range("K89").Value = 16238.75
$ws.Range("L89").Value = 13637.5
$ws.Range("M89").Value = -10622.75
$ws.Range("N89").Value = -24869.5
$ws.Range("H107").Value = 2091.4
$ws.Range("I107").Value = 1193
$ws.Range("K107").Value = 1193
$ws.Range("M107").Value = 727

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 774.4167
$ws.Range("I31").Value = 774.4167
$ws.Range("K31").Value = 774.4167
$ws.Range("M31").Value = -479.4167
$ws.Range("H34").Value = 774.4167
$ws.Range("I34").Value = 774.4167
$ws.Range("K34").Value = 774.4167
$ws.Range("M34").Value = -572.4167
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").Value = $null
$ws.Range("H86").Value = 6807.1
$ws.Range("I86").Value = 6600.8
$ws.Range("J86").Value = 7013.4
$ws.Range("K86").Value = 6600.8
$ws.Range("L86").Value = 7013.4
$ws.Range("M86").Value = -5477.8
$ws.Range("N86").Value = -9259.4
$ws.Range("H89").Value = 6807.1
$ws.Range("I89").Value = 6600.8
$ws.Range("J89").Value = 7013.4
$ws.Range("K89").Value = 33004
$ws.Range("L89").Value = 35067
$ws.Range("M89").Value = -27388
$ws.Range("N89").Value = -46299
$ws.Range("H134").Value = 4925.3335
$ws.Range("I134").Value = 4925.3335
$ws.Range("K134").Value = 14776.0005
$ws.Range("M134").Value = -12241.0005
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").Value = $null

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 400
$ws.Range("I18").Value = 400
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 1200
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -1031
$ws.Range("N18").Value = $null
$ws.Range("H23").Value = 315.5
$ws.Range("J23").Value = 320.66666
$ws.Range("L23").Value = 961.9999799999999
$ws.Range("N23").Value = -1431.99998
$ws.Range("H47").Value = 173.25
$ws.Range("I47").Value = 146.5
$ws.Range("J47").Value = 200
$ws.Range("K47").Value = 439.5
$ws.Range("L47").Value = 600
$ws.Range("M47").Value = -8.5
$ws.Range("N47").Value = -1462
$ws.Range("H113").Value = 846.8182
$ws.Range("I113").Value = 738.8570999999999
$ws.Range("J113").Value = 1035.75
$ws.Range("K113").Value = 2216.5713
$ws.Range("L113").Value = 3107.25
$ws.Range("M113").Value = -46.57129999999961
$ws.Range("N113").Value = -7447.25
$ws.Range("H123").Value = 5766.4707
$ws.Range("J123").Value = 6000
$ws.Range("L123").Value = 18000
$ws.Range("N123").Value = -22900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 34921
$ws.Range("J46").Value = 34921
$ws.Range("L46").Value = 34921
$ws.Range("N46").Value = -35233
$ws.Range("H102").Value = 2727.4285
$ws.Range("I102").Value = 2727.4285
$ws.Range("K102").Value = 2727.4285
$ws.Range("M102").Value = -1105.4285
$ws.Range("H126").Value = 4649.75
$ws.Range("I126").Value = 2866.6667
$ws.Range("K126").Value = 8600.000100000001
$ws.Range("M126").Value = -6130.000100000001
$ws.Range("H132").Value = 2669.25
$ws.Range("I132").Value = 1471.2
$ws.Range("K132").Value = 4413.6
$ws.Range("M132").Value = -1883.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1747
$ws.Range("I7").Value = 1799
$ws.Range("J7").Value = 1729.6666
$ws.Range("K7").Value = 1799
$ws.Range("L7").Value = 1729.6666
$ws.Range("M7").Value = -1687
$ws.Range("N7").Value = -1953.6666
$ws.Range("H16").Value = 3896.875
$ws.Range("I16").Value = 3876.4
$ws.Range("J16").Value = 3931
$ws.Range("K16").Value = 3876.4
$ws.Range("L16").Value = 3931
$ws.Range("M16").Value = -3706.4
$ws.Range("N16").Value = -4271
$ws.Range("H126").Value = 1747
$ws.Range("I126").Value = 1799
$ws.Range("J126").Value = 1729.6666
$ws.Range("K126").Value = 5397
$ws.Range("L126").Value = 5188.9998
$ws.Range("M126").Value = -2927
$ws.Range("N126").Value = -10128.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1614
$ws.Range("J126").Value = 3166.3333
$ws.Range("L126").Value = 9166.3333
$ws.Range("N126").Value = -14438.9999
